$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the missing sequence value in F5, and bump every value below it
# by one so the F-column stays a contiguous 0-based row index (previously
# it skipped a number at row 5).
for ($r = 70; $r -ge 6; $r--) {
    $cell = $ws.Cells.Item($r, 6)
    $old = $cell.Value()
    $cell.Value = $old + 1
}
$ws.Cells.Item(5, 6).Value = 4
